# Scheduled-runner update: refresh Market Board price/profit figures on the Leve
# profit-tracking sheets (H/I/J/K/L/M/N columns) to the latest pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 94
$ws.Range("H94").Value = 1372
$ws.Range("I94").Value = 1372
$ws.Range("K94").Value = 1372
$ws.Range("M94").Value = -921

# Row 113
$ws.Range("H113").Value = 14059.857
$ws.Range("I113").Value = 8749.75
$ws.Range("J113").Value = 21140
$ws.Range("K113").Value = 8749.75
$ws.Range("L113").Value = 21140
$ws.Range("M113").Value = -5495.75
$ws.Range("N113").Value = -27648

# Row 116
$ws.Range("H116").Value = 7603.391
$ws.Range("I116").Value = 6465.25
$ws.Range("J116").Value = 8845
$ws.Range("K116").Value = 6465.25
$ws.Range("L116").Value = 8845
$ws.Range("M116").Value = -3023.25
$ws.Range("N116").Value = -15729

# Row 127
$ws.Range("H127").Value = 8474.071
$ws.Range("J127").Value = 12687.25
$ws.Range("L127").Value = 38061.75
$ws.Range("N127").Value = -47981.75

# Row 129
$ws.Range("H129").Value = 203310
$ws.Range("I129").Value = 403997.2
$ws.Range("J129").Value = 2622.8
$ws.Range("K129").Value = 1211991.6
$ws.Range("L129").Value = 7868.400000000001
$ws.Range("M129").Value = -1206991.6
$ws.Range("N129").Value = -17868.4

# Row 131
$ws.Range("H131").Value = 9660.818
$ws.Range("I131").Value = 7783.625
$ws.Range("J131").Value = 14666.667
$ws.Range("K131").Value = 23350.875
$ws.Range("L131").Value = 44000.001
$ws.Range("M131").Value = -18310.875
$ws.Range("N131").Value = -54080.001

# Row 138
$ws.Range("H138").Value = 3208.4055
$ws.Range("I138").Value = 1734.0667
$ws.Range("K138").Value = 5202.2001
$ws.Range("M138").Value = -62.20010000000002

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4945.891
$ws.Range("I32").Value = 5369.3125
$ws.Range("K32").Value = 5369.3125
$ws.Range("M32").Value = -5082.3125

# Row 110
$ws.Range("H110").Value = 3981.5
$ws.Range("I110").Value = 3219.28
$ws.Range("K110").Value = 3219.28
$ws.Range("M110").Value = -1174.28

# Row 131
$ws.Range("H131").Value = 45731.066
$ws.Range("J131").Value = 45731.066
$ws.Range("L131").Value = 45731.066
$ws.Range("N131").Value = -55811.066

$ws = $wb.Worksheets.Item("BSM")
# Row 10
$ws.Range("H10").Value = 7500
$ws.Range("I10").Value = 5000
$ws.Range("J10").Value = 10000
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = -4860
$ws.Range("N10").Value = -10280

# Row 54
$ws.Range("H54").Value = 20000
$ws.Range("I54").Value = 20000
$ws.Range("K54").Value = 20000
$ws.Range("M54").Value = -19516

# Row 86
$ws.Range("H86").Value = 13890561
$ws.Range("I86").Value = 14707035
$ws.Range("J86").Value = 10499.5
$ws.Range("K86").Value = 14707035
$ws.Range("L86").Value = 10499.5
$ws.Range("M86").Value = -14705912
$ws.Range("N86").Value = -12745.5

# Row 89
$ws.Range("H89").Value = 13890561
$ws.Range("I89").Value = 14707035
$ws.Range("J89").Value = 10499.5
$ws.Range("K89").Value = 73535175
$ws.Range("L89").Value = 52497.5
$ws.Range("M89").Value = -73529559
$ws.Range("N89").Value = -63729.5

# Row 105
$ws.Range("H105").Value = 2033.6154
$ws.Range("I105").Value = 843.3333
$ws.Range("K105").Value = 843.3333
$ws.Range("M105").Value = 903.6667

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1946.5625
$ws.Range("I58").Value = 1599.5454
$ws.Range("K58").Value = 1599.5454
$ws.Range("M58").Value = -1396.5454

# Row 122
$ws.Range("H122").Value = 246858.38
$ws.Range("I122").Value = 343183.34
$ws.Range("K122").Value = 1029550.02
$ws.Range("M122").Value = -1027100.02

# Row 132
$ws.Range("H132").Value = 4294
$ws.Range("I132").Value = 2277.5715
$ws.Range("K132").Value = 6832.7145
$ws.Range("M132").Value = -4302.7145

# Row 136
$ws.Range("H136").Value = 1946.5625
$ws.Range("I136").Value = 1599.5454
$ws.Range("K136").Value = 4798.6362
$ws.Range("M136").Value = -2248.6362

$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 2852.889
$ws.Range("J129").Value = 2938
$ws.Range("L129").Value = 8814
$ws.Range("N129").Value = -18814

# Row 132
$ws.Range("H132").Value = 2904
$ws.Range("I132").Value = 1986.25
$ws.Range("K132").Value = 17876.25
$ws.Range("M132").Value = -15346.25

$ws = $wb.Worksheets.Item("GSM")
# Row 96
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

# Row 102
$ws.Range("H102").Value = 4885.5713
$ws.Range("I102").Value = 4140.727
$ws.Range("K102").Value = 4140.727
$ws.Range("M102").Value = -2518.727

# Row 113
$ws.Range("H113").Value = 6487.278
$ws.Range("I113").Value = 5097.8887
$ws.Range("K113").Value = 5097.8887
$ws.Range("M113").Value = -2927.8887

# Row 122
$ws.Range("H122").Value = 3965.75
$ws.Range("I122").Value = 2370.1428
$ws.Range("J122").Value = 6199.6
$ws.Range("K122").Value = 7110.428400000001
$ws.Range("L122").Value = 18598.8
$ws.Range("M122").Value = -4660.428400000001
$ws.Range("N122").Value = -23498.8

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2732.5
$ws.Range("I7").Value = 2732.5
$ws.Range("K7").Value = 2732.5
$ws.Range("M7").Value = -2620.5

# Row 22
$ws.Range("H22").Value = 10101761
$ws.Range("I22").Value = 18182082
$ws.Range("K22").Value = 18182082
$ws.Range("M22").Value = -18181787

# Row 27
$ws.Range("H27").Value = 10101761
$ws.Range("I27").Value = 18182082
$ws.Range("K27").Value = 18182082
$ws.Range("M27").Value = -18181975

# Row 69
$ws.Range("H69").Value = 120163
$ws.Range("J69").Value = 120163
$ws.Range("L69").Value = 120163
$ws.Range("N69").Value = -121785

# Row 72
$ws.Range("H72").Value = 120163
$ws.Range("J72").Value = 120163
$ws.Range("L72").Value = 360489
$ws.Range("N72").Value = -368601

# Row 82
$ws.Range("H82").Value = 76925410
$ws.Range("I82").Value = 100002390
$ws.Range("J82").Value = 2136
$ws.Range("K82").Value = 100002390
$ws.Range("L82").Value = 2136
$ws.Range("M82").Value = -100002029
$ws.Range("N82").Value = -2858

# Row 85
$ws.Range("H85").Value = 76925410
$ws.Range("I85").Value = 100002390
$ws.Range("J85").Value = 2136
$ws.Range("K85").Value = 100002390
$ws.Range("L85").Value = 2136
$ws.Range("M85").Value = -100001142
$ws.Range("N85").Value = -4632

# Row 94
$ws.Range("H94").Value = 51581.25
$ws.Range("J94").Value = 51581.25
$ws.Range("L94").Value = 51581.25
$ws.Range("N94").Value = -52933.25

# Row 126
$ws.Range("H126").Value = 2732.5
$ws.Range("I126").Value = 2732.5
$ws.Range("K126").Value = 8197.5
$ws.Range("M126").Value = -5727.5

# Row 132
$ws.Range("H132").Value = 2692.625
$ws.Range("I132").Value = 2381.3333
$ws.Range("J132").Value = 3092.8572
$ws.Range("K132").Value = 7143.999899999999
$ws.Range("L132").Value = 9278.5716
$ws.Range("M132").Value = -4613.999899999999
$ws.Range("N132").Value = -14338.5716

# Row 139
$ws.Range("H139").Value = 84998
$ws.Range("J139").Value = 80000
$ws.Range("L139").Value = 80000
$ws.Range("N139").Value = -90280

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 16657.889
$ws.Range("J4").Value = 16657.889
$ws.Range("L4").Value = 16657.889
$ws.Range("N4").Value = -16883.889

# Row 51
$ws.Range("H51").Value = 27135.727
$ws.Range("I51").Value = 15250
$ws.Range("J51").Value = 33927.57
$ws.Range("K51").Value = 15250
$ws.Range("L51").Value = 33927.57
$ws.Range("M51").Value = -14740
$ws.Range("N51").Value = -34947.57

# Row 63
$ws.Range("H63").Value = 42832.668
$ws.Range("J63").Value = 42832.668
$ws.Range("L63").Value = 42832.668
$ws.Range("N63").Value = -44080.668

# Row 66
$ws.Range("H66").Value = 42832.668
$ws.Range("J66").Value = 42832.668
$ws.Range("L66").Value = 128498.004
$ws.Range("N66").Value = -134738.004

# Row 95
$ws.Range("H95").Value = 31074.75
$ws.Range("J95").Value = 31074.75
$ws.Range("L95").Value = 31074.75
$ws.Range("N95").Value = -36566.75

# Row 107
$ws.Range("H107").Value = 402.3
$ws.Range("I107").Value = 402.3
$ws.Range("K107").Value = 1206.9
$ws.Range("M107").Value = 713.0999999999999

# Row 126
$ws.Range("H126").Value = 4066.3928
$ws.Range("J126").Value = 3576
$ws.Range("L126").Value = 10728
$ws.Range("N126").Value = -15668

# Row 135
$ws.Range("H135").Value = 50333.332
$ws.Range("J135").Value = 50333.332
$ws.Range("L135").Value = 50333.332
$ws.Range("N135").Value = -60473.332
